$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Cells.Item(2, 4).Value = '37.475.96'
$ws.Cells.Item(2, 5).Value = '  -0.05%  '
$ws.Cells.Item(3, 4).Value = '2.022.01'
$ws.Cells.Item(3, 5).Value = '  +0.45%  '
$ws.Cells.Item(4, 5).Value = '  -0.06%  '
$ws.Cells.Item(5, 4).NumberFormat = '@'
$ws.Cells.Item(5, 4).Value = '253.49'
$ws.Cells.Item(5, 5).Value = '  +2.79%  '
$ws.Cells.Item(6, 4).NumberFormat = '@'
$ws.Cells.Item(6, 4).Value = '0.618'
$ws.Cells.Item(6, 5).Value = '  -2.34%  '
$ws.Cells.Item(7, 5).Value = '  -0.03%  '
$ws.Cells.Item(8, 4).NumberFormat = '@'
$ws.Cells.Item(8, 4).Value = '56.56'
$ws.Cells.Item(8, 5).Value = '  -8.49%  '
$ws.Cells.Item(9, 4).NumberFormat = '@'
$ws.Cells.Item(9, 4).Value = '0.380'
$ws.Cells.Item(9, 5).Value = '  -1.39%  '
$ws.Cells.Item(10, 4).NumberFormat = '@'
$ws.Cells.Item(10, 4).Value = '0.0780'
$ws.Cells.Item(10, 5).Value = '  -3.29%  '
$ws.Cells.Item(11, 4).NumberFormat = '@'
$ws.Cells.Item(11, 4).Value = '0.101'
$ws.Cells.Item(11, 5).Value = '  -2.42%  '
$ws.Cells.Item(12, 5).Value = '  -3.51%  '
$ws.Cells.Item(13, 4).Value = '2.322.22'
$ws.Cells.Item(13, 5).Value = '  +0.88%  '
$ws.Cells.Item(14, 4).NumberFormat = '@'
$ws.Cells.Item(14, 4).Value = '0.810'
$ws.Cells.Item(14, 5).Value = '  -5.33%  '
$ws.Cells.Item(15, 4).NumberFormat = '@'
$ws.Cells.Item(15, 4).Value = '20.92'
$ws.Cells.Item(15, 5).Value = '  -9.24%  '
$ws.Cells.Item(16, 4).NumberFormat = '@'
$ws.Cells.Item(16, 4).Value = '5.31'
$ws.Cells.Item(16, 5).Value = '  -2.71%  '
$ws.Cells.Item(17, 4).Value = '2.032.74'
$ws.Cells.Item(17, 5).Value = '  +0.93%  '
$ws.Cells.Item(18, 4).Value = '37.373.63'
$ws.Cells.Item(18, 5).Value = '  -0.20%  '
$ws.Cells.Item(19, 4).NumberFormat = '@'
$ws.Cells.Item(19, 4).Value = '69.44'
$ws.Cells.Item(19, 5).Value = '  -1.75%  '
$ws.Cells.Item(20, 4).Value = '0.0₃0844'
$ws.Cells.Item(20, 5).Value = '  -2.99%  '
$ws.Cells.Item(21, 4).NumberFormat = '@'
$ws.Cells.Item(21, 4).Value = '5.17'
$ws.Cells.Item(21, 5).Value = '  -1.22%  '
$ws.Cells.Item(22, 4).NumberFormat = '@'
$ws.Cells.Item(22, 4).Value = '227.88'
$ws.Cells.Item(22, 5).Value = '  -1.89%  '
$ws.Cells.Item(23, 2).Value = 'Dai'
$ws.Cells.Item(23, 3).Value = 'https://coinranking.com/coin/MoTuySvg7+dai-dai'
$ws.Cells.Item(23, 4).NumberFormat = '@'
$ws.Cells.Item(23, 4).Value = '1.00'
$ws.Cells.Item(23, 5).Value = '  -0.03%  '
$ws.Cells.Item(24, 2).Value = 'PancakeSwap'
$ws.Cells.Item(24, 3).Value = 'https://coinranking.com/coin/ncYFcP709+pancakeswap-cake'
$ws.Cells.Item(24, 4).NumberFormat = '@'
$ws.Cells.Item(24, 4).Value = '2.60'
$ws.Cells.Item(24, 5).Value = '  +2.77%  '
$ws.Cells.Item(25, 4).NumberFormat = '@'
$ws.Cells.Item(25, 4).Value = '2.32'
$ws.Cells.Item(25, 5).Value = '  -2.35%  '
$ws.Cells.Item(26, 4).NumberFormat = '@'
$ws.Cells.Item(26, 4).Value = '163.19'
$ws.Cells.Item(26, 5).Value = '  -0.61%  '
$ws.Cells.Item(27, 4).NumberFormat = '@'
$ws.Cells.Item(27, 4).Value = '8.98'
$ws.Cells.Item(27, 5).Value = '  -4.55%  '
$ws.Cells.Item(28, 2).Value = 'EthereumClassic'
$ws.Cells.Item(28, 3).Value = 'https://coinranking.com/coin/hnfQfsYfeIGUQ+ethereumclassic-etc'
$ws.Cells.Item(28, 4).NumberFormat = '@'
$ws.Cells.Item(28, 4).Value = '19.73'
$ws.Cells.Item(28, 5).Value = '  -0.37%  '
$ws.Cells.Item(29, 2).Value = 'Kaspa'
$ws.Cells.Item(29, 3).Value = 'https://coinranking.com/coin/V8GxkwWow+kaspa-kas'
$ws.Cells.Item(29, 4).NumberFormat = '@'
$ws.Cells.Item(29, 4).Value = '0.131'
$ws.Cells.Item(29, 5).Value = '  -9.08%  '
$ws.Cells.Item(30, 5).Value = '  -0.31%  '
$ws.Cells.Item(31, 5).Value = '  -1.30%  '
$ws.Cells.Item(32, 4).NumberFormat = '@'
$ws.Cells.Item(32, 4).Value = '0.0667'
$ws.Cells.Item(32, 5).Value = '  +5.90%  '
$ws.Cells.Item(33, 5).Value = '  -4.43%  '
$ws.Cells.Item(34, 4).NumberFormat = '@'
$ws.Cells.Item(34, 4).Value = '4.53'
$ws.Cells.Item(34, 5).Value = '  -1.53%  '
$ws.Cells.Item(35, 5).Value = '  +1.92%  '
$ws.Cells.Item(36, 5).Value = '  -0.08%  '
$ws.Cells.Item(37, 5).Value = '  +0.21%  '
$ws.Cells.Item(38, 4).NumberFormat = '@'
$ws.Cells.Item(38, 4).Value = '3.39'
$ws.Cells.Item(38, 5).Value = '  -0.01%  '
$ws.Cells.Item(39, 4).NumberFormat = '@'
$ws.Cells.Item(39, 4).Value = '5.30'
$ws.Cells.Item(39, 5).Value = '  -5.08%  '
$ws.Cells.Item(40, 5).Value = '  +2.97%  '
$ws.Cells.Item(41, 4).NumberFormat = '@'
$ws.Cells.Item(41, 4).Value = '0.0960'
$ws.Cells.Item(41, 5).Value = '  -2.37%  '
$ws.Cells.Item(42, 5).Value = '  +1.23%  '
$ws.Cells.Item(43, 5).Value = '  -0.66%  '
$ws.Cells.Item(44, 4).Value = '1.411.85'
$ws.Cells.Item(44, 5).Value = '  +1.74%  '
$ws.Cells.Item(45, 2).Value = 'InjectiveProtocol'
$ws.Cells.Item(45, 3).Value = 'https://coinranking.com/coin/PkY9BmsyW+injectiveprotocol-inj'
$ws.Cells.Item(45, 4).NumberFormat = '@'
$ws.Cells.Item(45, 4).Value = '15.88'
$ws.Cells.Item(45, 5).Value = '  -6.06%  '
$ws.Cells.Item(46, 2).Value = 'Aave'
$ws.Cells.Item(46, 3).Value = 'https://coinranking.com/coin/ixgUfzmLR+aave-aave'
$ws.Cells.Item(46, 4).NumberFormat = '@'
$ws.Cells.Item(46, 4).Value = '90.29'
$ws.Cells.Item(46, 5).Value = '  -1.59%  '
$ws.Cells.Item(47, 5).Value = '  -1.95%  '
$ws.Cells.Item(48, 4).NumberFormat = '@'
$ws.Cells.Item(48, 4).Value = '7.28'
$ws.Cells.Item(48, 5).Value = '  -0.43%  '
$ws.Cells.Item(49, 5).Value = '  +0.21%  '
$ws.Cells.Item(50, 4).NumberFormat = '@'
$ws.Cells.Item(50, 4).Value = '2.01'
$ws.Cells.Item(50, 5).Value = '  +0.72%  '
$ws.Cells.Item(51, 4).Value = '2.214.23'
$ws.Cells.Item(51, 5).Value = '  +0.89%  '
